$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Loh")

$ws.Range("C11").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C6").Value = "Erste Besprechung und Plaung"
$ws.Range("D6").Value = 1.5

$ws.Range("C7").Value = "Einlesung in die Programmierung von Alexa und Installation des Programms für EV3"
$ws.Range("D6").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D7").Value = 4
$ws.Rows.Item(7).RowHeight = 45

$ws.Range("D23").Formula = "=SUM(D6:D18)"
